# ajout ssr, had et psy
# Add a new data row (row 46) to the "Feuil1" sheet:
#   B46 = "c"    (existing shared string, same as other rows in column B)
#   D46 = 182    (next position value, following row 45's 181)
#   F46 = "ZAD"  (new code name)
# Then move/leave the active selection on the newly added cell B46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B46").Value = "c"
$ws.Range("D46").Value = 182
$ws.Range("F46").Value = "ZAD"

$ws.Range("B46").Select() | Out-Null
